$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.953.52"
$ws.Range("E2").Value = "  +0.28%  "
$ws.Range("D3").Value = "1.641.03"
$ws.Range("E3").Value = "  +0.18%  "
$ws.Range("E4").Value = "  -0.07%  "
$r = $ws.Range("D5")
$r.NumberFormat = "@"
$r.Value = "215.29"
$r.Style = "Normal"
$ws.Range("E5").Value = "  -0.04%  "
$r = $ws.Range("D6")
$r.NumberFormat = "@"
$r.Value = "0.5079"
$r.Style = "Normal"
$ws.Range("E6").Value = "  +0.88%  "
$ws.Range("E7").Value = "  +0.11%  "
$r = $ws.Range("D8")
$r.NumberFormat = "@"
$r.Value = "0.2561"
$r.Style = "Normal"
$ws.Range("E8").Value = "  -0.05%  "
$r = $ws.Range("D9")
$r.NumberFormat = "@"
$r.Value = "0.06377"
$r.Style = "Normal"
$ws.Range("E9").Value = "  -0.18%  "
$r = $ws.Range("D10")
$r.NumberFormat = "@"
$r.Value = "19.49"
$r.Style = "Normal"
$ws.Range("E10").Value = "  -0.84%  "
$r = $ws.Range("D11")
$r.NumberFormat = "@"
$r.Value = "0.07774"
$r.Style = "Normal"
$ws.Range("E11").Value = "  +0.59%  "
$r = $ws.Range("D12")
$r.NumberFormat = "@"
$r.Value = "4.294"
$r.Style = "Normal"
$ws.Range("E12").Value = "  +0.83%  "
$ws.Range("D13").Value = "1.650.79"
$ws.Range("E13").Value = "  +0.87%  "
$r = $ws.Range("D14")
$r.NumberFormat = "@"
$r.Value = "0.5462"
$r.Style = "Normal"
$ws.Range("E14").Value = "  +0.36%  "
$ws.Range("E15").Value = "  -0.78%  "
$r = $ws.Range("D16")
$r.NumberFormat = "@"
$r.Value = "64.42"
$r.Style = "Normal"
$ws.Range("E16").Value = "  +0.24%  "
$ws.Range("D17").Value = "26.002.00"
$ws.Range("E17").Value = "  +0.41%  "
$r = $ws.Range("D18")
$r.NumberFormat = "@"
$r.Value = "1.003"
$r.Style = "Normal"
$ws.Range("E18").Value = "  -0.03%  "
$r = $ws.Range("D19")
$r.NumberFormat = "@"
$r.Value = "198.09"
$r.Style = "Normal"
$ws.Range("E19").Value = "  -2.37%  "
$r = $ws.Range("D20")
$r.NumberFormat = "@"
$r.Value = "4.439"
$r.Style = "Normal"
$ws.Range("E20").Value = "  +1.48%  "
$r = $ws.Range("D21")
$r.NumberFormat = "@"
$r.Value = "9.966"
$r.Style = "Normal"
$ws.Range("E21").Value = "  +0.71%  "
$r = $ws.Range("D22")
$r.NumberFormat = "@"
$r.Value = "6.059"
$r.Style = "Normal"
$ws.Range("E23").Value = "  +0.21%  "
$r = $ws.Range("D24")
$r.NumberFormat = "@"
$r.Value = "1.878"
$r.Style = "Normal"
$ws.Range("E24").Value = "  -2.50%  "
$r = $ws.Range("D25")
$r.NumberFormat = "@"
$r.Value = "140.97"
$r.Style = "Normal"
$r = $ws.Range("D26")
$r.NumberFormat = "@"
$r.Value = "0.1144"
$r.Style = "Normal"
$ws.Range("E26").Value = "  +0.58%  "
$r = $ws.Range("D27")
$r.NumberFormat = "@"
$r.Value = "6.880"
$r.Style = "Normal"
$ws.Range("E27").Value = "  +2.21%  "
$r = $ws.Range("D28")
$r.NumberFormat = "@"
$r.Value = "15.74"
$r.Style = "Normal"
$ws.Range("E28").Value = "  +0.50%  "
$r = $ws.Range("D29")
$r.NumberFormat = "@"
$r.Value = "1.238"
$r.Style = "Normal"
$ws.Range("E29").Value = "  -0.33%  "
$r = $ws.Range("D30")
$r.NumberFormat = "@"
$r.Value = "0.05032"
$r.Style = "Normal"
$ws.Range("E30").Value = "  +1.58%  "
$r = $ws.Range("D31")
$r.NumberFormat = "@"
$r.Value = "3.263"
$r.Style = "Normal"
$ws.Range("E31").Value = "  -0.29%  "
$r = $ws.Range("D32")
$r.NumberFormat = "@"
$r.Value = "3.188"
$r.Style = "Normal"
$ws.Range("E32").Value = "  +0.19%  "
$r = $ws.Range("D33")
$r.NumberFormat = "@"
$r.Value = "1.544"
$r.Style = "Normal"
$ws.Range("E33").Value = "  +0.17%  "
$r = $ws.Range("D34")
$r.NumberFormat = "@"
$r.Value = "2.367"
$r.Style = "Normal"
$ws.Range("E34").Value = "  -0.21%  "
$r = $ws.Range("D35")
$r.NumberFormat = "@"
$r.Value = "0.8993"
$r.Style = "Normal"
$ws.Range("E35").Value = "  +0.97%  "
$ws.Range("E36").Value = "  -1.26%  "
$ws.Range("D37").Value = "1.133.33"
$ws.Range("E37").Value = "  -2.26%  "
$r = $ws.Range("D38")
$r.NumberFormat = "@"
$r.Value = "0.5495"
$r.Style = "Normal"
$ws.Range("E38").Value = "  -1.82%  "
$ws.Range("E39").Value = "  +14.96%  "
$r = $ws.Range("D40")
$r.NumberFormat = "@"
$r.Value = "0.01557"
$r.Style = "Normal"
$ws.Range("E40").Value = "  -0.43%  "
$ws.Range("B41").Value = "mCoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin"
$r = $ws.Range("D41")
$r.NumberFormat = "@"
$r.Value = "2.549"
$r.Style = "Normal"
$ws.Range("E41").Value = "  -0.81%  "
$ws.Range("B42").Value = "PaxDollar"
$ws.Range("C42").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$r = $ws.Range("D42")
$r.NumberFormat = "@"
$r.Value = "1.003"
$r.Style = "Normal"
$ws.Range("E42").Value = "  +0.09%  "
$r = $ws.Range("D43")
$r.NumberFormat = "@"
$r.Value = "5.617"
$r.Style = "Normal"
$ws.Range("E43").Value = "  -0.54%  "
$r = $ws.Range("D44")
$r.NumberFormat = "@"
$r.Value = "0.8191"
$r.Style = "Normal"
$ws.Range("E44").Value = "  +1.57%  "
$r = $ws.Range("D45")
$r.NumberFormat = "@"
$r.Value = "100.16"
$r.Style = "Normal"
$ws.Range("E45").Value = "  +0.33%  "
$ws.Range("D46").Value = "1.779.56"
$ws.Range("E46").Value = "  +0.26%  "
$r = $ws.Range("D47")
$r.NumberFormat = "@"
$r.Value = "0.4530"
$r.Style = "Normal"
$ws.Range("E47").Value = "  -0.22%  "
$r = $ws.Range("D48")
$r.NumberFormat = "@"
$r.Value = "1.002"
$r.Style = "Normal"
$ws.Range("E48").Value = "  -0.01%  "
$r = $ws.Range("D49")
$r.NumberFormat = "@"
$r.Value = "54.91"
$r.Style = "Normal"
$ws.Range("E49").Value = "  +0.02%  "
$r = $ws.Range("D50")
$r.NumberFormat = "@"
$r.Value = "0.05073"
$r.Style = "Normal"
$ws.Range("E50").Value = "  +0.33%  "
$r = $ws.Range("D51")
$r.NumberFormat = "@"
$r.Value = "1.006"
$r.Style = "Normal"
$ws.Range("E51").Value = "  +0.42%  "
